# add update save feature
$wb = $excel.ActiveWorkbook

# --- "save" sheet: add an "id" column, turn the running total into a real number ---
$saveSheet = $wb.Worksheets.Item("save")

# new header cell F1 ("id") - clone the existing header formatting (bold,
# centered, thin-bordered) from the neighboring header cell, then set the text
$saveSheet.Range("E1").Copy($saveSheet.Range("F1"))
$saveSheet.Range("F1").Value = "id"

# the running total is now a real number instead of text, plus the new id value
$saveSheet.Range("E2").Value = 9000
$saveSheet.Range("F2").Value = 0

# --- "save_log" sheet: log individual save entries with an amount column ---
$logSheet = $wb.Worksheets.Item("save_log")

# new header cell C1 ("amount") - same formatting clone trick
$logSheet.Range("B1").Copy($logSheet.Range("C1"))
$logSheet.Range("C1").Value = "amount"

# row 2 becomes entry "a" with its amount
$logSheet.Range("B2").Value = "a"
$logSheet.Range("C2").Value = 6000

# row 3 is a new entry "b" with its amount
$logSheet.Range("A3").Value = 0
$logSheet.Range("B3").Value = "b"
$logSheet.Range("C3").Value = 3000

# the author was looking at the save_log tab when they saved
$logSheet.Activate()
